$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The source data had two rows in each of these pairs accidentally
# transposed. This swaps the full row content (columns B..AC; column A
# is the row's own running index and stays put) between each pair of
# rows to restore the correct order.
$rowPairs = @(
    @(129, 130),
    @(154, 156),
    @(157, 158),
    @(168, 169),
    @(210, 211)
)

foreach ($pair in $rowPairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    $rng1 = $ws.Range("B$r1" + ":AC$r1")
    $rng2 = $ws.Range("B$r2" + ":AC$r2")

    $vals1 = $rng1.Value2
    $vals2 = $rng2.Value2

    $rng1.Value2 = $vals2
    $rng2.Value2 = $vals1
}
